$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (K / Strike#) values per regenerated save_data
$ws.Range("G2").Value  = 0
$ws.Range("G3").Value  = 2
$ws.Range("G4").Value  = 2
$ws.Range("G5").Value  = 0
$ws.Range("G6").Value  = 3
$ws.Range("G7").Value  = 0
$ws.Range("G8").Value  = 2
$ws.Range("G9").Value  = 1
$ws.Range("G10").Value = 1
$ws.Range("G11").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("G14").Value = 2
$ws.Range("G16").Value = 0
